# Append new NAV history rows (2024-08-28 .. 2024-09-25) to the end of the
# existing data table, extending the used range from A1:J612 to A1:J633.
#
# Columns: A=Date, C..H=basket constituent values / basket value,
#          I=daily return, J=NAV. Column B only holds header/quantity rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRowsData = @"
613|2024-08-28|1890.599975585938|711.7999877929688|78.95999908447266|299.9500122070312|1323.25|29275.72984313965|0.0|338.796762044303
614|2024-08-29|1893.25|700.7999877929688|77.0|296.2000122070312|1298.650024414062|28928.15008544922|-0.01187262485180639|334.7743551875443
615|2024-08-30|1899.349975585938|696.0999755859375|75.83999633789062|299.2999877929688|1302.900024414062|28869.96923828125|-0.002011219071945895|334.1010506195927
616|2024-09-02|1885.400024414062|683.5999755859375|73.80999755859375|296.8999938964844|1303.849975585938|28487.57955932617|-0.01324524026329871|329.6758019319157
617|2024-09-03|1901.949951171875|689.4000244140625|74.47000122070312|297.1499938964844|1320.25|28715.86001586914|0.008013332830455756|332.3176038589435
618|2024-09-04|1901.300048828125|688.9500122070312|74.16000366210938|298.9500122070312|1327.099975585938|28721.03076171875|0.0001800658537390795|332.3774429119949
619|2024-09-05|1879.449951171875|687.5|76.0|290.6000061035156|1312.349975585938|28602.74978637695|-0.004118270556621158|331.0086226751654
620|2024-09-06|1872.349975585938|673.5499877929688|74.72000122070312|283.6000061035156|1289.699951171875|28191.60983276367|-0.01437414083205038|326.2506581162095
621|2024-09-09|1892.400024414062|664.1500244140625|74.33999633789062|281.5499877929688|1237.150024414062|28036.46997070312|-0.005503050836077006|324.4552841592924
622|2024-09-10|1922.449951171875|664.5999755859375|78.05000305175781|285.75|1250.300048828125|28561.49984741211|0.01872667555001102|330.5312529962301
623|2024-09-11|1957.599975585938|689.75|81.94999694824219|288.0499877929688|1237.699951171875|29297.64938354492|0.0257741904334731|339.0504284551695
624|2024-09-12|1996.400024414062|729.1500244140625|81.69999694824219|291.7000122070312|1237.300048828125|29912.80038452148|0.02099659917843318|346.1693344027187
625|2024-09-13|1988.050048828125|713.7000122070312|83.11000061035156|289.9500122070312|1241.5|29812.18057250977|-0.003363771051799782|345.004900016634
626|2024-09-16|1989.900024414062|714.2000122070312|84.69999694824219|290.3999938964844|1226.599975585938|29926.49987792969|0.003834650911961043|346.3278733711138
627|2024-09-17|2006.550048828125|731.0999755859375|82.0|284.2999877929688|1193.800048828125|29823.24993896484|-0.003450117433913109|345.1330015373461
628|2024-09-18|1987.800048828125|743.25|80.81999969482422|282.8500061035156|1166.400024414062|29685.31034851074|-0.004625236710834788|343.536679708515
629|2024-09-19|1998.599975585938|735.9500122070312|80.97000122070312|272.7000122070312|1121.300048828125|29455.8603515625|-0.0077294120982553|340.8813431401815
630|2024-09-20|2048.10009765625|746.5|83.44999694824219|277.3500061035156|1149.400024414062|30118.95037841797|0.02251131078642199|348.5550289969032
631|2024-09-23|2082.39990234375|773.9500122070312|82.88999938964844|286.2999877929688|1162.75|30664.31942749023|0.01810717313253569|354.8663752531661
632|2024-09-24|2068.14990234375|781.8499755859375|83.79000091552734|291.7999877929688|1141.199951171875|30770.6690826416|0.003468188994144963|356.0971189102113
633|2024-09-25|2061.60009765625|775.8499755859375|82.95999908447266|289.8500061035156|1118.449951171875|30542.33015441895|-0.007420668286718119|353.4546403129226
"@

$rows = $newRowsData.Trim() -split "`r?`n"

foreach ($line in $rows) {
    $f = $line.Trim() -split '\|'

    $r  = [int]$f[0]
    $dt = $f[1]

    # Leading apostrophe forces the date-looking string to be stored as
    # literal text, matching the existing column A cells (e.g. "2022-05-11")
    # instead of being auto-parsed into a date serial number. ClearFormats
    # drops the transient "quote prefix" cell style so the cell ends up with
    # the same default (unstyled) formatting as its neighbours.
    $ws.Cells.Item($r, 1).Value = "'" + $dt
    $ws.Cells.Item($r, 1).ClearFormats()

    $ws.Cells.Item($r, 3).Value = [double]$f[2]   # C
    $ws.Cells.Item($r, 4).Value = [double]$f[3]   # D
    $ws.Cells.Item($r, 5).Value = [double]$f[4]   # E
    $ws.Cells.Item($r, 6).Value = [double]$f[5]   # F
    $ws.Cells.Item($r, 7).Value = [double]$f[6]   # G
    $ws.Cells.Item($r, 8).Value = [double]$f[7]   # H
    $ws.Cells.Item($r, 9).Value = [double]$f[8]   # I
    $ws.Cells.Item($r, 10).Value = [double]$f[9]  # J
}
